$d = $word.ActiveDocument

# 1. Remove the old auto "_GoBack" bookmark first (Word keeps only a single
#    instance, tracking the most recent edit location; it will be re-created
#    at the new edit location below).
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# 2. Insert the new list paragraph "Die Designer unterschätzen..." after the
#    "Die Designer müssen..." item, before the blank spacing paragraph (numId 32 list).
$p24 = $d.Paragraphs(24)
$target = $d.Range($p24.Range.Start, $p24.Range.End)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="32"/></w:numPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="MS Reference Sans Serif" w:hAnsi="MS Reference Sans Serif" w:cs="Arial"/><w:szCs w:val="28"/><w:lang w:val="de-AT" w:eastAsia="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="MS Reference Sans Serif" w:hAnsi="MS Reference Sans Serif" w:cs="Arial"/><w:szCs w:val="28"/><w:lang w:val="de-AT" w:eastAsia="en-GB"/></w:rPr><w:t>Die Designer untersch&#228;tzen die Komplexit&#228;t der Umsetzung der Wireframes</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MS Reference Sans Serif" w:hAnsi="MS Reference Sans Serif" w:cs="Arial"/><w:szCs w:val="28"/><w:lang w:val="de-AT" w:eastAsia="en-GB"/></w:rPr><w:t>.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="MS Reference Sans Serif" w:hAnsi="MS Reference Sans Serif" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="28"/><w:lang w:val="de-AT" w:eastAsia="en-GB"/></w:rPr></w:pPr></w:p>'
$target.InsertXML($xml)

Write-Output "done"
